$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: SBY | 22 | 40 | 2024
$ws.Range("A3").Value = "SBY"
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("C2").Copy($ws.Range("C3"))

$ws.Range("ZZ1").Formula = '="2024"'
$ws.Range("ZZ1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()

# Row 4: ZSY | 22 | 30 | 2024
$ws.Range("A4").Value = "ZSY"
$ws.Range("B2").Copy($ws.Range("B4"))

$ws.Range("ZZ1").Formula = '="30"'
$ws.Range("ZZ1").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()

$ws.Range("D3").Copy($ws.Range("D4"))
